# Added slides and recording for Abel Brodeur workshop
#
# - "2024 - Fall" sheet, row 5 (Sept. 25th / Abel Broduer - Reproducibility,
#   Replication Packets, and Pre-analysis Plans workshop): record in-person
#   (G5) and Zoom (H5) attendance, which ripples through the running totals,
#   the summary stats and the "Attendance Descriptives" helper columns.
# - A "JMC" note is added to the Notes column (J) on several rows across the
#   "2024 - Fall", "2024 - Spring" and "2023 - Fall" sheets.
# - A few leftover selection / frozen-pane view settings are nudged to match
#   where the author last clicked while doing this edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "2024 - Fall": Abel Brodeur workshop attendance + Notes column updates
# ---------------------------------------------------------------------
$fall2024 = $wb.Worksheets.Item("2024 - Fall")
$fall2024.Activate()

# Attendance for the Abel Broduer workshop (row 5)
$fall2024.Range("G5").Value = 5
$fall2024.Range("H5").Value = 5

# "JMC" notes
$fall2024.Range("J2").Value = "JMC"
$fall2024.Range("J6").Value = "JMC"
$fall2024.Range("J10").Value = "JMC"
$fall2024.Range("J11").Value = "JMC"
$fall2024.Range("J12").Value = "JMC"
$fall2024.Range("J16").Value = "JMC"

# Notes column got narrower once the shorter "JMC" notes were added
$fall2024.Columns.Item(10).ColumnWidth = 25.83

$fall2024.Range("J6").Select() | Out-Null

# ---------------------------------------------------------------------
# "Attendance Descriptives": just a leftover selection change
# ---------------------------------------------------------------------
$descriptives = $wb.Worksheets.Item("Attendance Descriptives")
$descriptives.Activate()
$descriptives.Range("O26").Select() | Out-Null

# ---------------------------------------------------------------------
# "2024 - Spring": add "JMC" notes
# ---------------------------------------------------------------------
$spring2024 = $wb.Worksheets.Item("2024 - Spring")
$spring2024.Activate()

$spring2024.Range("J6").Value = "JMC"
$spring2024.Range("J10").Value = "JMC"
$spring2024.Range("J14").Value = "JMC"

$spring2024.Range("J16").Select() | Out-Null

# ---------------------------------------------------------------------
# "2023 - Fall": add "JMC" notes
# ---------------------------------------------------------------------
$fall2023 = $wb.Worksheets.Item("2023 - Fall")
$fall2023.Activate()

$fall2023.Range("J2").Value = "JMC"
$fall2023.Range("J4").Value = "JMC"
$fall2023.Range("J11").Value = "JMC"
$fall2023.Range("J14").Value = "JMC"
$fall2023.Range("J15").Value = "JMC"

$fall2023.Range("A7").Select() | Out-Null
$fall2023.Range("J16").Select() | Out-Null

# Leave the workbook back on the tab that was active before ("2024 - Fall")
$fall2024.Activate()
$fall2024.Range("J6").Select() | Out-Null
